$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format column D and E as Text first so numeric-looking strings (e.g. "18.20", "1.00")
# are preserved exactly as typed, matching the source data which is stored as text.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "67.423.50"
$ws.Range("E2").Value = "  -3.74%  "

# Row 3
$ws.Range("D3").Value = "3.782.22"
$ws.Range("E3").Value = "  -3.79%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "596.91"
$ws.Range("E5").Value = "  -2.18%  "

# Row 6
$ws.Range("D6").Value = "166.59"
$ws.Range("E6").Value = "  -2.48%  "

# Row 7
$ws.Range("D7").Value = "3.786.40"
$ws.Range("E7").Value = "  -3.77%  "

# Row 8
$ws.Range("E8").Value = "  -0.07%  "

# Row 9
$ws.Range("D9").Value = "0.523"
$ws.Range("E9").Value = "  -2.75%  "

# Row 10
$ws.Range("D10").Value = "0.162"
$ws.Range("E10").Value = "  -4.82%  "

# Row 11
$ws.Range("E11").Value = "  +0.07%  "

# Row 12
$ws.Range("D12").Value = "0.453"
$ws.Range("E12").Value = "  -3.59%  "

# Row 13
$ws.Range("D13").Value = "0.0000257"
$ws.Range("E13").Value = "  -0.30%  "

# Row 14
$ws.Range("D14").Value = "36.61"
$ws.Range("E14").Value = "  -4.97%  "

# Row 15
$ws.Range("D15").Value = "4.427.04"
$ws.Range("E15").Value = "  -3.64%  "

# Row 16
$ws.Range("D16").Value = "3.812.14"
$ws.Range("E16").Value = "  -3.41%  "

# Row 17
$ws.Range("D17").Value = "67.513.27"
$ws.Range("E17").Value = "  -3.64%  "

# Row 18
$ws.Range("D18").Value = "18.20"
$ws.Range("E18").Value = "  -2.86%  "

# Row 19
$ws.Range("D19").Value = "7.33"
$ws.Range("E19").Value = "  -4.26%  "

# Row 20
$ws.Range("E20").Value = "  -1.29%  "

# Row 21
$ws.Range("D21").Value = "10.87"
$ws.Range("E21").Value = "  -2.56%  "

# Row 22
$ws.Range("D22").Value = "463.44"
$ws.Range("E22").Value = "  -6.30%  "

# Row 23
$ws.Range("D23").Value = "0.725"
$ws.Range("E23").Value = "  -3.34%  "

# Row 24
$ws.Range("D24").Value = "0.0000159"
$ws.Range("E24").Value = "  -5.01%  "

# Row 25
$ws.Range("D25").Value = "82.42"
$ws.Range("E25").Value = "  -4.31%  "

# Row 26
$ws.Range("D26").Value = "2.20"
$ws.Range("E26").Value = "  -3.90%  "

# Row 27
$ws.Range("D27").Value = "12.00"
$ws.Range("E27").Value = "  -2.73%  "

# Row 28
$ws.Range("D28").Value = "0.997"
$ws.Range("E28").Value = "  -0.30%  "

# Row 29
$ws.Range("D29").Value = "9.96"
$ws.Range("E29").Value = "  -2.16%  "

# Row 30
$ws.Range("D30").Value = "2.94"
$ws.Range("E30").Value = "  -2.17%  "

# Row 31
$ws.Range("D31").Value = "3.942.33"
$ws.Range("E31").Value = "  -3.49%  "

# Row 32
$ws.Range("D32").Value = "7.56"
$ws.Range("E32").Value = "  -4.07%  "

# Row 33
$ws.Range("D33").Value = "31.07"
$ws.Range("E33").Value = "  -3.74%  "

# Row 34
$ws.Range("D34").Value = "2.28"
$ws.Range("E34").Value = "  -7.16%  "

# Row 35
$ws.Range("D35").Value = "9.41"
$ws.Range("E35").Value = "  -2.30%  "

# Row 36
$ws.Range("D36").Value = "3.756.33"
$ws.Range("E36").Value = "  -3.58%  "

# Row 37
$ws.Range("D37").Value = "0.103"
$ws.Range("E37").Value = "  -4.91%  "

# Row 38
$ws.Range("D38").Value = "3.62"
$ws.Range("E38").Value = "  +9.39%  "

# Row 39
$ws.Range("B39").Value = "Mantle"
$ws.Range("C39").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D39").Value = "1.01"
$ws.Range("E39").Value = "  -3.73%  "

# Row 40
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "0.139"
$ws.Range("E40").Value = "  -2.27%  "

# Row 41
$ws.Range("D41").Value = "5.86"
$ws.Range("E41").Value = "  -4.96%  "

# Row 42
$ws.Range("E42").Value = "  +0.47%  "

# Row 43
$ws.Range("D43").Value = "0.311"
$ws.Range("E43").Value = "  -5.89%  "

# Row 44
$ws.Range("D44").Value = "1.96"
$ws.Range("E44").Value = "  -8.13%  "

# Row 45
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.01%  "

# Row 46
$ws.Range("B46").Value = "FLOKI"
$ws.Range("C46").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D46").Value = "0.000293"
$ws.Range("E46").Value = "  +5.24%  "

# Row 47
$ws.Range("D47").Value = "8.64"
$ws.Range("E47").Value = "  -0.55%  "

# Row 48
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").Value = "414.26"
$ws.Range("E48").Value = "  -5.68%  "

# Row 49
$ws.Range("D49").Value = "46.66"
$ws.Range("E49").Value = "  -3.62%  "

# Row 50
$ws.Range("D50").Value = "141.69"
$ws.Range("E50").Value = "  -1.12%  "

# Row 51
$ws.Range("D51").Value = "26.02"
$ws.Range("E51").Value = "  +2.69%  "

